$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header label in H1 (results for an "Explicit mode" column)
$ws.Range("H1").Value = "Explicit mode"

# Match the authored column width for the new column H
$ws.Columns.Item(8).ColumnWidth = 11

# Leave the selection where the author left it after typing the new header
$ws.Range("H3").Select()
